# Insert a new weekly price record as row 74 in the "Arveja Verde" sheet.
# This pushes the existing rows 74-76 down to 75-77 (dimension grows to A1:R77).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 74, shifting rows 74-76 down to 75-77.
$ws.Rows.Item(74).Insert()

# Populate the new row 74 with the added weekly record.
$ws.Cells.Item(74, 1).Value  = 9
$ws.Cells.Item(74, 2).Value  = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(74, 3).Value  = "Metropolitana"
$ws.Cells.Item(74, 4).Value  = 44516
$ws.Cells.Item(74, 5).Value  = 13
$ws.Cells.Item(74, 6).Value  = 100112022
$ws.Cells.Item(74, 7).Value  = "Arveja Verde"
$ws.Cells.Item(74, 8).Value  = "Sin especificar"
$ws.Cells.Item(74, 9).Value  = "Primera"
$ws.Cells.Item(74, 10).Value = 34
$ws.Cells.Item(74, 11).Value = 14000
$ws.Cells.Item(74, 12).Value = 15000
$ws.Cells.Item(74, 13).Value = 14500
$ws.Cells.Item(74, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(74, 15).Value = "Región Metropolitana"
$ws.Cells.Item(74, 16).Value = 580
$ws.Cells.Item(74, 17).Value = 25
$ws.Cells.Item(74, 18).Value = "Hortaliza"
